# Replace the stale "congenital" category label with "misc_long_term"
# across every worksheet in the workbook that still uses it.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Value2 -eq "congenital") {
            $cell.Value = "misc_long_term"
        }
    }
}
